$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 18984.61550902891
$ws.Range("C2").Value = 21187.1530531127
$ws.Range("D2").Value = 32953.54879383396
